$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G columns for rows 334-362 (AgTests / AgPosit) per the commit diff
$updates = @{
    334 = @{ F = 205310; G = 3470 }
    335 = @{ F = 130356; G = 2979 }
    336 = @{ F = 103758; G = 3296 }
    337 = @{ F = 106418; G = 2987 }
    338 = @{ F = 222698; G = 3118 }
    339 = @{ F = 652308; G = 5607 }
    341 = @{ F = 296517; G = 3650 }
    342 = @{ F = 179181; G = 3044 }
    343 = @{ F = 131312; G = 2934 }
    344 = @{ F = 135174; G = 2502 }
    345 = @{ F = 287536; G = 3290 }
    346 = @{ F = 660518; G = 4766 }
    347 = @{ F = 333216 }
    348 = @{ F = 232581; G = 3230 }
    349 = @{ F = 158709; G = 2745 }
    350 = @{ F = 127881; G = 2770 }
    351 = @{ F = 148374; G = 2795 }
    352 = @{ F = 302455; G = 3505 }
    353 = @{ F = 708188; G = 5173 }
    354 = @{ F = 303415; G = 2767 }
    355 = @{ F = 219137; G = 3380 }
    356 = @{ F = 158356; G = 2857 }
    357 = @{ F = 136797; G = 2999 }
    358 = @{ F = 159021; G = 2652 }
    359 = @{ F = 315852; G = 3321 }
    360 = @{ F = 714769; G = 4816 }
    361 = @{ F = 324919; G = 2554 }
    362 = @{ F = 216669; G = 3034 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Append new row 363 with the new day's data
$ws.Range("A363").Value = 44257
$ws.Range("A363").NumberFormat = "yyyy-mm-dd"
$ws.Range("B363").Value = 0
$ws.Range("C363").Value = -2060990
$ws.Range("D363").Value = -311002
$ws.Range("E363").Value = 7489
$ws.Range("F363").Value = 154772
$ws.Range("G363").Value = 2817

$wb.Save()
